$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be written as TEXT,
# matching the source inlineStr cells, without leaving a NumberFormat
# behind on the cell (Style reset to Normal afterwards).
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = '66.095.35'
$ws.Cells.Item(2, 5).Value = '  +1.60%  '

$ws.Cells.Item(3, 4).Value = '3.213.62'
$ws.Cells.Item(3, 5).Value = '  +1.04%  '

$ws.Cells.Item(4, 5).Value = '  +0.05%  '

Set-TextValue $ws.Cells.Item(5, 4) '604.50'
$ws.Cells.Item(5, 5).Value = '  +4.23%  '

Set-TextValue $ws.Cells.Item(6, 4) '153.62'
$ws.Cells.Item(6, 5).Value = '  +0.89%  '

$ws.Cells.Item(7, 5).Value = '  +0.08%  '

$ws.Cells.Item(8, 4).Value = '3.211.43'
$ws.Cells.Item(8, 5).Value = '  +1.05%  '

$ws.Cells.Item(9, 5).Value = '  +0.33%  '

$ws.Cells.Item(10, 5).Value = '  -1.19%  '

$ws.Cells.Item(11, 5).Value = '  -1.33%  '

$ws.Cells.Item(12, 5).Value = '  +0.98%  '

Set-TextValue $ws.Cells.Item(13, 4) '0.0000272'
$ws.Cells.Item(13, 5).Value = '  -0.19%  '

Set-TextValue $ws.Cells.Item(14, 4) '38.59'
$ws.Cells.Item(14, 5).Value = '  +1.70%  '

$ws.Cells.Item(15, 4).Value = '3.741.55'
$ws.Cells.Item(15, 5).Value = '  +1.04%  '

$ws.Cells.Item(16, 2).Value = 'Polkadot'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Cells.Item(16, 4) '7.47'
$ws.Cells.Item(16, 5).Value = '  +3.96%  '

$ws.Cells.Item(17, 2).Value = 'WrappedBTC'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(17, 4).Value = '66.209.08'
$ws.Cells.Item(17, 5).Value = '  +1.64%  '

$ws.Cells.Item(18, 4).Value = '3.208.95'
$ws.Cells.Item(18, 5).Value = '  +0.83%  '

$ws.Cells.Item(19, 5).Value = '  -0.24%  '

Set-TextValue $ws.Cells.Item(20, 4) '511.60'
$ws.Cells.Item(20, 5).Value = '  -0.62%  '

Set-TextValue $ws.Cells.Item(21, 4) '15.53'
$ws.Cells.Item(21, 5).Value = '  +4.17%  '

Set-TextValue $ws.Cells.Item(22, 4) '0.733'
$ws.Cells.Item(22, 5).Value = '  +0.18%  '

Set-TextValue $ws.Cells.Item(23, 4) '15.27'
$ws.Cells.Item(23, 5).Value = '  +0.05%  '

Set-TextValue $ws.Cells.Item(24, 4) '8.02'
$ws.Cells.Item(24, 5).Value = '  +2.55%  '

Set-TextValue $ws.Cells.Item(25, 4) '85.19'
$ws.Cells.Item(25, 5).Value = '  -0.35%  '

$ws.Cells.Item(27, 5).Value = '  +2.49%  '

Set-TextValue $ws.Cells.Item(28, 4) '9.18'
$ws.Cells.Item(28, 5).Value = '  +1.95%  '

Set-TextValue $ws.Cells.Item(29, 4) '2.25'
$ws.Cells.Item(29, 5).Value = '  +2.79%  '

Set-TextValue $ws.Cells.Item(30, 4) '2.86'
$ws.Cells.Item(30, 5).Value = '  +2.97%  '

$ws.Cells.Item(31, 2).Value = 'EthereumClassic'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Cells.Item(31, 4) '28.17'
$ws.Cells.Item(31, 5).Value = '  +0.70%  '

$ws.Cells.Item(32, 2).Value = 'NEARProtocol'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Cells.Item(32, 4) '6.80'
$ws.Cells.Item(32, 5).Value = '  +7.01%  '

Set-TextValue $ws.Cells.Item(33, 4) '1.21'
$ws.Cells.Item(33, 5).Value = '  +0.62%  '

$ws.Cells.Item(34, 5).Value = '  +0.14%  '

Set-TextValue $ws.Cells.Item(35, 4) '6.61'
$ws.Cells.Item(35, 5).Value = '  +0.21%  '

Set-TextValue $ws.Cells.Item(36, 4) '55.35'
$ws.Cells.Item(36, 5).Value = '  -0.77%  '

Set-TextValue $ws.Cells.Item(37, 4) '0.0906'

Set-TextValue $ws.Cells.Item(38, 4) '480.81'
$ws.Cells.Item(38, 5).Value = '  +0.89%  '

Set-TextValue $ws.Cells.Item(39, 4) '0.0421'
$ws.Cells.Item(39, 5).Value = '  -0.58%  '

Set-TextValue $ws.Cells.Item(40, 4) '2.97'
$ws.Cells.Item(40, 5).Value = '  -6.34%  '

Set-TextValue $ws.Cells.Item(41, 4) '8.85'
$ws.Cells.Item(41, 5).Value = '  +1.90%  '

Set-TextValue $ws.Cells.Item(42, 4) '0.297'
$ws.Cells.Item(42, 5).Value = '  +3.49%  '

$ws.Cells.Item(43, 5).Value = '  -0.51%  '

$ws.Cells.Item(44, 4).Value = '2.942.69'
$ws.Cells.Item(44, 5).Value = '  -4.32%  '

Set-TextValue $ws.Cells.Item(45, 4) '2.45'
$ws.Cells.Item(45, 5).Value = '  +1.89%  '

$ws.Cells.Item(46, 4).Value = '0.0₃0641'
$ws.Cells.Item(46, 5).Value = '  +3.99%  '

Set-TextValue $ws.Cells.Item(47, 4) '28.86'
$ws.Cells.Item(47, 5).Value = '  -1.19%  '

$ws.Cells.Item(48, 5).Value = '  +0.01%  '

$ws.Cells.Item(49, 5).Value = '  -0.16%  '

Set-TextValue $ws.Cells.Item(50, 4) '2.31'
$ws.Cells.Item(50, 5).Value = '  +2.01%  '

Set-TextValue $ws.Cells.Item(51, 4) '33.87'
$ws.Cells.Item(51, 5).Value = '  +3.94%  '
